$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-04 Thursday" "2024-04-05 Friday"

Replace-Text "718÷8=" "115÷7="
Replace-Text "127÷5=" "563÷9="
Replace-Text "630÷6=" "269÷5="
Replace-Text "446÷9=" "794÷6="
Replace-Text "801÷6=" "184÷6="

Replace-Text "251÷6=" "113÷3="
Replace-Text "481÷8=" "902÷2="
Replace-Text "943÷7=" "553÷8="
Replace-Text "899÷9=" "690÷3="
Replace-Text "323÷6=" "474÷7="

Replace-Text "592÷8=" "228÷2="
Replace-Text "590÷7=" "189÷4="
Replace-Text "377÷7=" "713÷5="
Replace-Text "386÷3=" "490÷7="
Replace-Text "967÷2=" "526÷3="

Replace-Text "901÷7=" "281÷9="
Replace-Text "900÷5=" "296÷8="
Replace-Text "291÷6=" "505÷4="
Replace-Text "530÷2=" "876÷2="
Replace-Text "470÷2=" "588÷3="

Replace-Text "783÷3=" "765÷8="
Replace-Text "959÷6=" "653÷6="
Replace-Text "308÷8=" "628÷6="
Replace-Text "649÷2=" "276÷9="
Replace-Text "281÷8=" "285÷8="
